# Auto-generated edit script: refresh market-price derived values in Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 488.2
$ws.Range("I33").Value = 467.6154
$ws.Range("J33").Value = 526.4286
$ws.Range("K33").Value = 467.6154
$ws.Range("L33").Value = 526.4286
$ws.Range("M33").Value = -238.6154
$ws.Range("N33").Value = -984.4286
$ws.Range("H117").Value = 48267.332
$ws.Range("J117").Value = 48267.332
$ws.Range("L117").Value = 48267.332
$ws.Range("N117").Value = -57445.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2940.9048
$ws.Range("I63").Value = 2663.0625
$ws.Range("K63").Value = 2663.0625
$ws.Range("M63").Value = -1977.0625
$ws.Range("H66").Value = 2940.9048
$ws.Range("I66").Value = 2663.0625
$ws.Range("K66").Value = 13315.3125
$ws.Range("M66").Value = -9883.3125
$ws.Range("H80").Value = 51401.855
$ws.Range("J80").Value = 51401.855
$ws.Range("L80").Value = 51401.855
$ws.Range("N80").Value = -53397.855
$ws.Range("H83").Value = 51401.855
$ws.Range("J83").Value = 51401.855
$ws.Range("L83").Value = 154205.565
$ws.Range("N83").Value = -164189.565
$ws.Range("H104").Value = 40723.332
$ws.Range("J104").Value = 40723.332
$ws.Range("L104").Value = 40723.332
$ws.Range("N104").Value = -47711.332
$ws.Range("H105").Value = 47960
$ws.Range("J105").Value = 47960
$ws.Range("L105").Value = 47960
$ws.Range("N105").Value = -54948
$ws.Range("H107").Value = 36260
$ws.Range("J107").Value = 36260
$ws.Range("L107").Value = 36260
$ws.Range("N107").Value = -43940
$ws.Range("H109").Value = 43152.8
$ws.Range("J109").Value = 43152.8
$ws.Range("L109").Value = 43152.8
$ws.Range("N109").Value = -45926.8
$ws.Range("H111").Value = 45471
$ws.Range("J111").Value = 45471
$ws.Range("L111").Value = 45471
$ws.Range("N111").Value = -53651
$ws.Range("H113").Value = 46336.668
$ws.Range("J113").Value = 46336.668
$ws.Range("L113").Value = 46336.668
$ws.Range("N113").Value = -55014.668
$ws.Range("H114").Value = 40365
$ws.Range("J114").Value = 40365
$ws.Range("L114").Value = 40365
$ws.Range("N114").Value = -49043
$ws.Range("H117").Value = 47924.168
$ws.Range("J117").Value = 47924.168
$ws.Range("L117").Value = 47924.168
$ws.Range("N117").Value = -57102.168
$ws.Range("H118").Value = 49372.332
$ws.Range("J118").Value = 49372.332
$ws.Range("L118").Value = 49372.332
$ws.Range("N118").Value = -52686.332
$ws.Range("H119").Value = 52592.5
$ws.Range("J119").Value = 52592.5
$ws.Range("L119").Value = 52592.5
$ws.Range("N119").Value = -62268.5
$ws.Range("H123").Value = 40878
$ws.Range("J123").Value = 40878
$ws.Range("L123").Value = 40878
$ws.Range("N123").Value = -50678
$ws.Range("H131").Value = 40297.168
$ws.Range("J131").Value = 40297.168
$ws.Range("L131").Value = 40297.168
$ws.Range("N131").Value = -50377.168

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 48694
$ws.Range("J110").Value = 48694
$ws.Range("L110").Value = 48694
$ws.Range("N110").Value = -56874
$ws.Range("H111").Value = 47694
$ws.Range("J111").Value = 47694
$ws.Range("L111").Value = 47694
$ws.Range("N111").Value = -55874
$ws.Range("H112").Value = 45939.5
$ws.Range("J112").Value = 45939.5
$ws.Range("L112").Value = 45939.5
$ws.Range("N112").Value = -48893.5
$ws.Range("H117").Value = 49914
$ws.Range("J117").Value = 49914
$ws.Range("L117").Value = 49914
$ws.Range("N117").Value = -59092
$ws.Range("H119").Value = 47992
$ws.Range("J119").Value = 47992
$ws.Range("L119").Value = 47992
$ws.Range("N119").Value = -57668
$ws.Range("H120").Value = 48761
$ws.Range("J120").Value = 48761
$ws.Range("L120").Value = 48761
$ws.Range("N120").Value = -58437
$ws.Range("H126").Value = 42441.332
$ws.Range("J126").Value = 42441.332
$ws.Range("L126").Value = 42441.332
$ws.Range("N126").Value = -52321.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 202089.5
$ws.Range("I31").Value = 2945.889
$ws.Range("J31").Value = 230997.44
$ws.Range("K31").Value = 2945.889
$ws.Range("L31").Value = 230997.44
$ws.Range("M31").Value = -2650.889
$ws.Range("N31").Value = -231587.44
$ws.Range("H34").Value = 202089.5
$ws.Range("I34").Value = 2945.889
$ws.Range("J34").Value = 230997.44
$ws.Range("K34").Value = 2945.889
$ws.Range("L34").Value = 230997.44
$ws.Range("M34").Value = -2743.889
$ws.Range("N34").Value = -231401.44
$ws.Range("H62").Value = 3028.5
$ws.Range("J62").Value = 4003
$ws.Range("L62").Value = 4003
$ws.Range("N62").Value = -5251
$ws.Range("H65").Value = 3028.5
$ws.Range("J65").Value = 4003
$ws.Range("L65").Value = 20015
$ws.Range("N65").Value = -26255
$ws.Range("H109").Value = 29165.555
$ws.Range("J109").Value = 29165.555
$ws.Range("L109").Value = 29165.555
$ws.Range("N109").Value = -31245.555
$ws.Range("H111").Value = 47267
$ws.Range("J111").Value = 47267
$ws.Range("L111").Value = 47267
$ws.Range("N111").Value = -55447
$ws.Range("H116").Value = 47793.5
$ws.Range("J116").Value = 47793.5
$ws.Range("L116").Value = 47793.5
$ws.Range("N116").Value = -56971.5
$ws.Range("H118").Value = 44734
$ws.Range("J118").Value = 44734
$ws.Range("L118").Value = 44734
$ws.Range("N118").Value = -48048
$ws.Range("H119").Value = 48250.668
$ws.Range("J119").Value = 48250.668
$ws.Range("L119").Value = 48250.668
$ws.Range("N119").Value = -57926.668
$ws.Range("H120").Value = 31823.75
$ws.Range("J120").Value = 31823.75
$ws.Range("L120").Value = 31823.75
$ws.Range("N120").Value = -39081.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 55567350
$ws.Range("I137").Value = 3970
$ws.Range("J137").Value = 83349040
$ws.Range("K137").Value = 11910
$ws.Range("L137").Value = 250047120
$ws.Range("M137").Value = -6810
$ws.Range("N137").Value = -250057320

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 44947.5
$ws.Range("J104").Value = 44947.5
$ws.Range("L104").Value = 44947.5
$ws.Range("N104").Value = -51935.5
$ws.Range("H114").Value = 41130.168
$ws.Range("J114").Value = 41130.168
$ws.Range("L114").Value = 41130.168
$ws.Range("N114").Value = -49808.168
$ws.Range("H116").Value = 38939.285
$ws.Range("J116").Value = 38939.285
$ws.Range("L116").Value = 38939.285
$ws.Range("N116").Value = -48117.285
$ws.Range("H118").Value = 38201.332
$ws.Range("J118").Value = 38201.332
$ws.Range("L118").Value = 38201.332
$ws.Range("N118").Value = -41515.332
$ws.Range("H130").Value = 44339.8
$ws.Range("J130").Value = 44339.8
$ws.Range("L130").Value = 44339.8
$ws.Range("N130").Value = -54379.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 35277
$ws.Range("J109").Value = 35277
$ws.Range("L109").Value = 35277
$ws.Range("N109").Value = -38051
$ws.Range("H110").Value = 45553
$ws.Range("J110").Value = 45553
$ws.Range("L110").Value = 45553
$ws.Range("N110").Value = -53733
$ws.Range("H111").Value = 43886
$ws.Range("J111").Value = 43886
$ws.Range("L111").Value = 43886
$ws.Range("N111").Value = -52066
$ws.Range("H112").Value = 35206.8
$ws.Range("J112").Value = 35206.8
$ws.Range("L112").Value = 35206.8
$ws.Range("N112").Value = -38160.8
$ws.Range("H116").Value = 45668
$ws.Range("J116").Value = 45668
$ws.Range("L116").Value = 45668
$ws.Range("N116").Value = -54846
$ws.Range("H120").Value = 51737.5
$ws.Range("J120").Value = 51737.5
$ws.Range("L120").Value = 51737.5
$ws.Range("N120").Value = -61413.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 50058.5
$ws.Range("J105").Value = 50058.5
$ws.Range("L105").Value = 50058.5
$ws.Range("N105").Value = -57046.5
$ws.Range("H110").Value = 23096
$ws.Range("J110").Value = 23096
$ws.Range("L110").Value = 23096
$ws.Range("N110").Value = -31276
$ws.Range("H116").Value = 48836
$ws.Range("J116").Value = 48836
$ws.Range("L116").Value = 48836
$ws.Range("N116").Value = -58014
$ws.Range("H117").Value = 47332
$ws.Range("J117").Value = 47332
$ws.Range("L117").Value = 47332
$ws.Range("N117").Value = -56510
$ws.Range("H118").Value = 42384
$ws.Range("J118").Value = 42384
$ws.Range("L118").Value = 42384
$ws.Range("N118").Value = -45698
